$d = $word.ActiveDocument

# 1. Replace the first paragraph's text. The original paragraph has four
#    runs: the base sentence plus three colored "(This is a change..."
#    runs. Collapse them into a single plain run with the final text.
$p1 = $d.Paragraphs.Item(1)
$start = $p1.Range.Start
$end = $p1.Range.End - 1   # exclude the paragraph mark
$r = $d.Range($start, $end)
$r.Text = "This is a Microsoft word document."

# 2. Remove the trailing empty shaded paragraph (the last paragraph in the
#    body, just before the sectPr).
$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIndex)
$pPrev = $d.Paragraphs.Item($lastIndex - 1)
$delStart = $pPrev.Range.End - 1
$delEnd = $pLast.Range.End
$d.Range($delStart, $delEnd).Delete()
